$d = $word.ActiveDocument

# 1. Apply strikethrough formatting to the first block of paragraphs that
#    were marked "done": Count donors .. Calculate the total amount of
#    donations received during a campaign (paragraphs 10-13).
foreach ($i in 10,11,12,13) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Font.StrikeThrough = 1
}

# 2. Apply strikethrough formatting to the second block: Search by Name,
#    First Name, Last name, Organization/Company Name (paragraphs 15-18).
foreach ($i in 15,16,17,18) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Font.StrikeThrough = 1
}

# 3. Remove the text "No" from the sub-bullet under "Search donors by
#    solicitation status" (paragraph 34), leaving the (now empty) bullet
#    in place.
$pNo = $d.Paragraphs.Item(34)
$rNo = $pNo.Range
$rNoTrim = $d.Range($rNo.Start, $rNo.End - 1)
$rNoTrim.Text = ""

# 4. Insert a new sub-bullet "Grant Name" right after
#    "Organization/Company Name" (paragraph 18), matching its
#    (now strikethrough) formatting.
$pOrg = $d.Paragraphs.Item(18)
$pOrg.Range.InsertParagraphAfter()
$d2 = $word.ActiveDocument
$pGrant = $d2.Paragraphs.Item(19)
$pGrant.Range.Text = "Grant Name"
